$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Egypt Division 1")

# Row 74
$ws.Range("B74").Value = 7217624
$ws.Range("F74").Value = 'Pyramids FC'
$ws.Range("G74").Value = 'Enppi'
$ws.Range("H74").Value = 1
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 1.444
$ws.Range("L74").Value = 3.75
$ws.Range("M74").Value = 6.5
$ws.Range("N74").Value = 1.5
$ws.Range("O74").Value = 3.75
$ws.Range("P74").Value = 6
$ws.Range("Q74").Value = -1
$ws.Range("R74").Value = 1.85
$ws.Range("S74").Value = 1.95
$ws.Range("U74").Value = 1.975
$ws.Range("V74").Value = 1.825
$ws.Range("W74").Value = 0.5
$ws.Range("Z74").Value = 0
$ws.Range("AA74").Value = -0
$ws.Range("AB74").Value = -1
$ws.Range("AC74").Value = 0.825

# Row 75
$ws.Range("B75").Value = 7217625
$ws.Range("F75").Value = 'Al Ahly Cairo'
$ws.Range("G75").Value = 'Ismaily SC'
$ws.Range("H75").Value = 3
$ws.Range("I75").Value = 1
$ws.Range("K75").Value = 1.25
$ws.Range("L75").Value = 4.75
$ws.Range("M75").Value = 10
$ws.Range("N75").Value = 1.333
$ws.Range("O75").Value = 4.2
$ws.Range("P75").Value = 8
$ws.Range("Q75").Value = -1.25
$ws.Range("R75").Value = 1.775
$ws.Range("S75").Value = 2.025
$ws.Range("U75").Value = 1.875
$ws.Range("V75").Value = 1.925
$ws.Range("W75").Value = 0.333
$ws.Range("Z75").Value = 0.7749999999999999
$ws.Range("AA75").Value = -1
$ws.Range("AB75").Value = 0.875
$ws.Range("AC75").Value = -1

# Row 79
$ws.Range("B79").Value = 7217629
$ws.Range("F79").Value = 'Enppi'
$ws.Range("G79").Value = 'National Bank'
$ws.Range("H79").Value = 3
$ws.Range("K79").Value = 2.4
$ws.Range("L79").Value = 2.875
$ws.Range("M79").Value = 3
$ws.Range("N79").Value = 2.5
$ws.Range("O79").Value = 2.8
$ws.Range("P79").Value = 2.9
$ws.Range("Q79").Value = 0
$ws.Range("R79").Value = 1.8
$ws.Range("S79").Value = 2
$ws.Range("T79").Value = 2
$ws.Range("U79").Value = 1.85
$ws.Range("V79").Value = 1.95
$ws.Range("W79").Value = 1.5
$ws.Range("Z79").Value = 0.8
$ws.Range("AB79").Value = 0.8500000000000001

# Row 80
$ws.Range("B80").Value = 7217630
$ws.Range("F80").Value = 'El Zamalek'
$ws.Range("G80").Value = 'Smouha'
$ws.Range("H80").Value = 5
$ws.Range("K80").Value = 1.65
$ws.Range("L80").Value = 3.5
$ws.Range("M80").Value = 4.5
$ws.Range("N80").Value = 1.909
$ws.Range("O80").Value = 3.25
$ws.Range("P80").Value = 3.6
$ws.Range("Q80").Value = -0.5
$ws.Range("R80").Value = 1.925
$ws.Range("S80").Value = 1.875
$ws.Range("T80").Value = 2.5
$ws.Range("U80").Value = 1.975
$ws.Range("V80").Value = 1.825
$ws.Range("W80").Value = 0.909
$ws.Range("Z80").Value = 0.925
$ws.Range("AB80").Value = 0.9750000000000001

# Row 87
$ws.Range("B87").Value = 7217638
$ws.Range("F87").Value = 'Al Ittihad Al Sakandary'
$ws.Range("G87").Value = 'Al Moqawloon Al Arab'
$ws.Range("H87").Value = 3
$ws.Range("I87").Value = 2
$ws.Range("J87").Value = 'H'
$ws.Range("K87").Value = 2.25
$ws.Range("L87").Value = 2.9
$ws.Range("M87").Value = 3.1
$ws.Range("N87").Value = 2.45
$ws.Range("O87").Value = 2.875
$ws.Range("P87").Value = 2.875
$ws.Range("Q87").Value = 0
$ws.Range("R87").Value = 1.725
$ws.Range("S87").Value = 2.075
$ws.Range("T87").Value = 2
$ws.Range("U87").Value = 1.775
$ws.Range("V87").Value = 2.025
$ws.Range("W87").Value = 1.45
$ws.Range("X87").Value = -1
$ws.Range("Z87").Value = 0.7250000000000001
$ws.Range("AA87").Value = -1
$ws.Range("AB87").Value = 0.7749999999999999
$ws.Range("AC87").Value = -1

# Row 88
$ws.Range("B88").Value = 7217639
$ws.Range("F88").Value = 'ZED FC'
$ws.Range("G88").Value = 'Smouha'
$ws.Range("H88").Value = 1
$ws.Range("I88").Value = 1
$ws.Range("J88").Value = 'D'
$ws.Range("K88").Value = 2.1
$ws.Range("L88").Value = 3
$ws.Range("M88").Value = 3.25
$ws.Range("N88").Value = 1.833
$ws.Range("O88").Value = 3.2
$ws.Range("P88").Value = 4
$ws.Range("Q88").Value = -0.5
$ws.Range("R88").Value = 1.9
$ws.Range("S88").Value = 1.9
$ws.Range("T88").Value = 2.25
$ws.Range("U88").Value = 2
$ws.Range("V88").Value = 1.8
$ws.Range("W88").Value = -1
$ws.Range("X88").Value = 2.2
$ws.Range("Z88").Value = -1
$ws.Range("AA88").Value = 0.8999999999999999
$ws.Range("AB88").Value = -0.5
$ws.Range("AC88").Value = 0.4

# Row 107
$ws.Range("B107").Value = 7217658
$ws.Range("F107").Value = 'Talaea El Geish'
$ws.Range("G107").Value = 'El Gounah'
$ws.Range("K107").Value = 2.1
$ws.Range("L107").Value = 2.875
$ws.Range("M107").Value = 3.5
$ws.Range("N107").Value = 2.15
$ws.Range("O107").Value = 2.8
$ws.Range("P107").Value = 3.5
$ws.Range("Q107").Value = -0.25
$ws.Range("R107").Value = 1.85
$ws.Range("S107").Value = 1.95
$ws.Range("U107").Value = 1.825
$ws.Range("V107").Value = 1.975
$ws.Range("X107").Value = 1.8
$ws.Range("Z107").Value = -0.5
$ws.Range("AA107").Value = 0.475

# Row 108
$ws.Range("B108").Value = 7217659
$ws.Range("F108").Value = 'Pharco FC'
$ws.Range("G108").Value = 'Al Moqawloon Al Arab'
$ws.Range("K108").Value = 2.5
$ws.Range("L108").Value = 2.8
$ws.Range("M108").Value = 2.8
$ws.Range("N108").Value = 2.6
$ws.Range("O108").Value = 2.7
$ws.Range("P108").Value = 2.75
$ws.Range("Q108").Value = 0
$ws.Range("R108").Value = 1.825
$ws.Range("S108").Value = 1.975
$ws.Range("U108").Value = 2
$ws.Range("V108").Value = 1.8
$ws.Range("X108").Value = 1.7
$ws.Range("Z108").Value = 0
$ws.Range("AA108").Value = -0

# Row 122
$ws.Range("B122").Value = 7549591
$ws.Range("F122").Value = 'Pharco FC'
$ws.Range("G122").Value = 'ZED FC'
$ws.Range("K122").Value = 2.875
$ws.Range("L122").Value = 2.75
$ws.Range("M122").Value = 2.5
$ws.Range("N122").Value = 2.9
$ws.Range("P122").Value = 2.45
$ws.Range("Q122").Value = 0
$ws.Range("R122").Value = 2.075
$ws.Range("S122").Value = 1.725
$ws.Range("T122").Value = 2
$ws.Range("U122").Value = 2
$ws.Range("V122").Value = 1.8
$ws.Range("Y122").Value = 1.45
$ws.Range("AA122").Value = 0.7250000000000001
$ws.Range("AB122").Value = 1

# Row 123
$ws.Range("B123").Value = 7549592
$ws.Range("F123").Value = 'Talaea El Geish'
$ws.Range("G123").Value = 'Al Moqawloon Al Arab'
$ws.Range("K123").Value = 2.375
$ws.Range("L123").Value = 2.8
$ws.Range("M123").Value = 3.3
$ws.Range("N123").Value = 2.25
$ws.Range("P123").Value = 3.5
$ws.Range("Q123").Value = -0.25
$ws.Range("R123").Value = 1.9
$ws.Range("S123").Value = 1.9
$ws.Range("T123").Value = 1.75
$ws.Range("U123").Value = 1.775
$ws.Range("V123").Value = 2.025
$ws.Range("Y123").Value = 2.5
$ws.Range("AA123").Value = 0.8999999999999999
$ws.Range("AB123").Value = 0.7749999999999999

# Row 148
$ws.Range("B148").Value = 7217715
$ws.Range("F148").Value = 'El Daklyeh'
$ws.Range("G148").Value = 'Pharco FC'
$ws.Range("H148").Value = 1
$ws.Range("I148").Value = 1
$ws.Range("J148").Value = 'D'
$ws.Range("K148").Value = 3
$ws.Range("L148").Value = 2.8
$ws.Range("M148").Value = 2.5
$ws.Range("N148").Value = 3.1
$ws.Range("O148").Value = 2.7
$ws.Range("P148").Value = 2.5
$ws.Range("Q148").Value = 0.25
$ws.Range("R148").Value = 1.75
$ws.Range("S148").Value = 2.05
$ws.Range("T148").Value = 1.75
$ws.Range("U148").Value = 1.925
$ws.Range("V148").Value = 1.875
$ws.Range("W148").Value = -1
$ws.Range("X148").Value = 1.7
$ws.Range("Z148").Value = 0.375
$ws.Range("AA148").Value = -0.5
$ws.Range("AB148").Value = 0.4625
$ws.Range("AC148").Value = -0.5

# Row 149
$ws.Range("B149").Value = 7217716
$ws.Range("F149").Value = 'El Masry'
$ws.Range("G149").Value = 'Baladiyet El Mahallah'
$ws.Range("H149").Value = 3
$ws.Range("I149").Value = 0
$ws.Range("J149").Value = 'H'
$ws.Range("K149").Value = 1.615
$ws.Range("L149").Value = 3.4
$ws.Range("M149").Value = 5.75
$ws.Range("N149").Value = 1.65
$ws.Range("O149").Value = 3.4
$ws.Range("P149").Value = 5.25
$ws.Range("Q149").Value = -0.75
$ws.Range("R149").Value = 1.825
$ws.Range("S149").Value = 1.975
$ws.Range("T149").Value = 2.25
$ws.Range("U149").Value = 1.9
$ws.Range("V149").Value = 1.9
$ws.Range("W149").Value = 0.6499999999999999
$ws.Range("X149").Value = -1
$ws.Range("Z149").Value = 0.825
$ws.Range("AA149").Value = -1
$ws.Range("AB149").Value = 0.8999999999999999
$ws.Range("AC149").Value = -1

# Row 173
$ws.Range("B173").Value = 7878949
$ws.Range("F173").Value = 'El Gounah'
$ws.Range("G173").Value = 'El Zamalek'
$ws.Range("H173").Value = 3
$ws.Range("I173").Value = 2
$ws.Range("K173").Value = 1.909
$ws.Range("L173").Value = 3.4
$ws.Range("M173").Value = 3.4
$ws.Range("N173").Value = 3.4
$ws.Range("O173").Value = 3.25
$ws.Range("P173").Value = 1.95
$ws.Range("Q173").Value = 0.5
$ws.Range("U173").Value = 1.95
$ws.Range("V173").Value = 1.85
$ws.Range("W173").Value = 2.4
$ws.Range("AB173").Value = 0.95

# Row 174
$ws.Range("B174").Value = 7878664
$ws.Range("F174").Value = 'National Bank'
$ws.Range("G174").Value = 'Al Ahly Cairo'
$ws.Range("H174").Value = 4
$ws.Range("I174").Value = 3
$ws.Range("K174").Value = 5
$ws.Range("L174").Value = 4
$ws.Range("M174").Value = 1.5
$ws.Range("N174").Value = 5.5
$ws.Range("O174").Value = 3.4
$ws.Range("P174").Value = 1.55
$ws.Range("Q174").Value = 1
$ws.Range("U174").Value = 1.8
$ws.Range("V174").Value = 2
$ws.Range("W174").Value = 4.5
$ws.Range("AB174").Value = 0.8

# Row 175
$ws.Range("U175").Value = 1.975
$ws.Range("V175").Value = 1.825

# Row 176
$ws.Range("B176").Value = 7217744
$ws.Range("F176").Value = 'El Gounah'
$ws.Range("G176").Value = 'Al Moqawloon Al Arab'
$ws.Range("K176").Value = 2.45
$ws.Range("M176").Value = 2.9
$ws.Range("N176").Value = 2.45
$ws.Range("P176").Value = 2.9
$ws.Range("R176").Value = 1.775
$ws.Range("S176").Value = 2.025

# Row 177
$ws.Range("B177").Value = 7217745
$ws.Range("F177").Value = 'Ismaily SC'
$ws.Range("G177").Value = 'Smouha'
$ws.Range("K177").Value = 2.7
$ws.Range("M177").Value = 2.6
$ws.Range("N177").Value = 2.7
$ws.Range("P177").Value = 2.6
$ws.Range("R177").Value = 2
$ws.Range("S177").Value = 1.8
$ws.Range("U177").Value = 1.825
$ws.Range("V177").Value = 1.975
